# Poster (fa) slide 13: shorten the "extended" deadline from 15 Khordad to
# 14 Khordad (May 5 -> May 4) and nudge the red "minus" strike-through shape
# that sits under the old date so it still spans the new, shorter text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

# The date text lives in "TextBox 11", nested inside the "Group 8" group.
$grp = $s.Shapes.Item("Group 8")
$tb = $grp.GroupItems.Item("TextBox 11")
$tr = $tb.TextFrame.TextRange

# 1) "31 " / "اردیبهشت " / "1400" (three separate runs) collapse into a
#    single run "31 اردیبهشت 1400".
$full = $tr.Text
$idx1 = $full.IndexOf("31 ")
$len1 = "31 اردیبهشت 1400".Length
$tr.Characters($idx1 + 1, $len1).Text = "31 اردیبهشت 1400"

# 2) "15 " becomes "14 " inside "تا 15 خرداد 1400 تمدید شد".
$full = $tr.Text
$idx2 = $full.IndexOf("15 ")
$tr.Characters($idx2 + 1, 3).Text = "14 "

# 3) The red math-minus strike shape is resized/repositioned to match the
#    new (shorter) date text. PowerPoint COM reports Left/Top/Width/Height
#    in points; the OOXML stores EMU (1 pt = 12700 EMU).
$minus = $s.Shapes.Item("Minus 13")
$minus.Left = 2524125 / 12700.0
$minus.Top = 5981700 / 12700.0
$minus.Width = 1812925 / 12700.0
$minus.Height = 94616 / 12700.0
